$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-6 (sending cluster reindexed; TPM-based stats recalculated) ---
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Calca"
$ws.Range("C2").Value = "Calcrl"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.21212
$ws.Range("H2").Value = 0.63636
$ws.Range("I2").Value = 0.2045603692733198
$ws.Range("J2").Value = 0.2783689285053439
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 40.688829
$ws.Range("N2").Value = 122.066487
$ws.Range("O2").Value = 0.5299680863154126
$ws.Range("P2").Value = 0.5381766555421038
$ws.Range("Q2").Value = 8.63091440748
$ws.Range("R2").Value = 77.67822966732
$ws.Range("S2").Value = 0.1084104674397554
$ws.Range("T2").Value = 0.149811658949845

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Calcrl"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.21212
$ws.Range("H3").Value = 0.63636
$ws.Range("I3").Value = 0.2045603692733198
$ws.Range("J3").Value = 0.2783689285053439
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.22188666666667
$ws.Range("N3").Value = 42.66566
$ws.Range("O3").Value = 0.1852387066860051
$ws.Range("P3").Value = 0.1881078318023236
$ws.Range("Q3").Value = 3.016746599733334
$ws.Range("R3").Value = 27.1507193976
$ws.Range("S3").Value = 0.03789249824340137
$ws.Range("T3").Value = 0.05236337558227627

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Calca"
$ws.Range("C4").Value = "Calcrl"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.21212
$ws.Range("H4").Value = 0.63636
$ws.Range("I4").Value = 0.2045603692733198
$ws.Range("J4").Value = 0.2783689285053439
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.24048933333333
$ws.Range("N4").Value = 30.721468
$ws.Range("O4").Value = 0.1333813891503258
$ws.Range("P4").Value = 0.1354473066926532
$ws.Range("Q4").Value = 2.172212597386667
$ws.Range("R4").Value = 19.54991337648
$ws.Range("S4").Value = 0.02728454621877902
$ws.Range("T4").Value = 0.03770432163296857

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Calca"
$ws.Range("C5").Value = "Calcrl"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.21212
$ws.Range("H5").Value = 0.63636
$ws.Range("I5").Value = 0.2045603692733198
$ws.Range("J5").Value = 0.2783689285053439
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.5130905
$ws.Range("N5").Value = 7.026181
$ws.Range("O5").Value = 0.04575766604976163
$ws.Range("P5").Value = 0.03097759823147425
$ws.Range("Q5").Value = 0.7451967568600001
$ws.Range("R5").Value = 4.471180541160001
$ws.Range("S5").Value = 0.009360205064224487
$ws.Range("T5").Value = 0.008623200827364523

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Calca"
$ws.Range("C6").Value = "Calcrl"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.21212
$ws.Range("H6").Value = 0.63636
$ws.Range("I6").Value = 0.2045603692733198
$ws.Range("J6").Value = 0.2783689285053439
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 8.111703
$ws.Range("N6").Value = 24.335109
$ws.Range("O6").Value = 0.1056541517984947
$ws.Range("P6").Value = 0.1072906077314452
$ws.Range("Q6").Value = 1.72065444036
$ws.Range("R6").Value = 15.48588996324
$ws.Range("S6").Value = 0.02161265230715945
$ws.Range("T6").Value = 0.02986637151288956

# --- Add new rows 7-11 (MuSCs as sending cluster) ---
# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Calca"
$ws.Range("C7").Value = "Calcrl"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8248355
$ws.Range("H7").Value = 1.649671
$ws.Range("I7").Value = 0.7954396307266801
$ws.Range("J7").Value = 0.721631071494656
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 40.688829
$ws.Range("N7").Value = 122.066487
$ws.Range("O7").Value = 0.5299680863154126
$ws.Range("P7").Value = 0.5381766555421038
$ws.Range("Q7").Value = 33.5615906126295
$ws.Range("R7").Value = 201.369543675777
$ws.Range("S7").Value = 0.4215576188756571
$ws.Range("T7").Value = 0.3883649965922588

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Calca"
$ws.Range("C8").Value = "Calcrl"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8248355
$ws.Range("H8").Value = 1.649671
$ws.Range("I8").Value = 0.7954396307266801
$ws.Range("J8").Value = 0.721631071494656
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 14.22188666666667
$ws.Range("N8").Value = 42.66566
$ws.Range("O8").Value = 0.1852387066860051
$ws.Range("P8").Value = 0.1881078318023236
$ws.Range("Q8").Value = 11.73071699964333
$ws.Range("R8").Value = 70.38430199786
$ws.Range("S8").Value = 0.1473462084426037
$ws.Range("T8").Value = 0.1357444562200473

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Calca"
$ws.Range("C9").Value = "Calcrl"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8248355
$ws.Range("H9").Value = 1.649671
$ws.Range("I9").Value = 0.7954396307266801
$ws.Range("J9").Value = 0.721631071494656
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.24048933333333
$ws.Range("N9").Value = 30.721468
$ws.Range("O9").Value = 0.1333813891503258
$ws.Range("P9").Value = 0.1354473066926532
$ws.Range("Q9").Value = 8.446719139504669
$ws.Range("R9").Value = 50.68031483702801
$ws.Range("S9").Value = 0.1060968429315468
$ws.Range("T9").Value = 0.09774298505968462

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Calca"
$ws.Range("C10").Value = "Calcrl"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8248355
$ws.Range("H10").Value = 1.649671
$ws.Range("I10").Value = 0.7954396307266801
$ws.Range("J10").Value = 0.721631071494656
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.5130905
$ws.Range("N10").Value = 7.026181
$ws.Range("O10").Value = 0.04575766604976163
$ws.Range("P10").Value = 0.03097759823147425
$ws.Range("Q10").Value = 2.89772175911275
$ws.Range("R10").Value = 11.590887036451
$ws.Range("S10").Value = 0.03639746098553714
$ws.Range("T10").Value = 0.02235439740410972

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Calca"
$ws.Range("C11").Value = "Calcrl"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8248355
$ws.Range("H11").Value = 1.649671
$ws.Range("I11").Value = 0.7954396307266801
$ws.Range("J11").Value = 0.721631071494656
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.111703
$ws.Range("N11").Value = 24.335109
$ws.Range("O11").Value = 0.1056541517984947
$ws.Range("P11").Value = 0.1072906077314452
$ws.Range("Q11").Value = 6.690820599856501
$ws.Range("R11").Value = 40.144923599139
$ws.Range("S11").Value = 0.0840414994913352
$ws.Range("T11").Value = 0.07742423621855558

